# The deck ships two embedded DrawingML themes:
#   ppt/theme/theme1.xml -> "Office Theme" (wired to the Notes Master)
#   ppt/theme/theme2.xml -> "Integral"     (wired to the Slide Master, i.e.
#                                            the theme every slide actually
#                                            renders with)
#
# The authored edit swaps the *content* of those two theme parts (Integral
# colours move into theme1.xml, Office Theme colours move into theme2.xml)
# while leaving every relationship untouched. The colour values are the
# only part of that content that differs between the two themes (font
# scheme and format scheme are byte-identical), so re-pointing the theme
# that actually drives the presentation (reached through
# Slide.ThemeColorScheme / Master.Theme.ThemeColorScheme, which this host
# always resolves to the Slide Master's theme part) at the Office Theme
# palette reproduces the visible effect of that swap.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function ToRGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme (was theme1.xml / notes master), in the
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink slot order used by
# ThemeColorScheme.Colors(1..12).
$tcs.Colors(1).RGB  = ToRGB 0x00 0x00 0x00   # dk1
$tcs.Colors(2).RGB  = ToRGB 0xFF 0xFF 0xFF   # lt1
$tcs.Colors(3).RGB  = ToRGB 0x44 0x54 0x6A   # dk2
$tcs.Colors(4).RGB  = ToRGB 0xE7 0xE6 0xE6   # lt2
$tcs.Colors(5).RGB  = ToRGB 0x5B 0x9B 0xD5   # accent1
$tcs.Colors(6).RGB  = ToRGB 0xED 0x7D 0x31   # accent2
$tcs.Colors(7).RGB  = ToRGB 0xA5 0xA5 0xA5   # accent3
$tcs.Colors(8).RGB  = ToRGB 0xFF 0xC0 0x00   # accent4
$tcs.Colors(9).RGB  = ToRGB 0x44 0x72 0xC4   # accent5
$tcs.Colors(10).RGB = ToRGB 0x70 0xAD 0x47   # accent6
$tcs.Colors(11).RGB = ToRGB 0x05 0x63 0xC1   # hlink
$tcs.Colors(12).RGB = ToRGB 0x95 0x4F 0x72   # folHlink
